# Fruta / hortaliza, semanal
# Insert a new weekly price-report row before row 14 (this pushes the
# existing rows 14-39 down to 15-40, which is exactly what the diff shows:
# every following row's content becomes the content of the row that used
# to precede it, and the sheet's used range grows from A1:T39 to A1:T40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 14. Excel copies the
# formatting (incl. the date number format on column D) from the row
# above automatically, same as the interactive "Insert" command.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with this week's record.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44925
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101001
$ws.Range("J14").Value = "Arándano (blue)"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 400
$ws.Range("N14").Value = 3000
$ws.Range("O14").Value = 3500
$ws.Range("P14").Value = 3250
$ws.Range("Q14").Value = "$/bandeja 2 kilos"
$ws.Range("R14").Value = "Provincia de Curicó"
$ws.Range("S14").Value = 1625
$ws.Range("T14").Value = 2
